# Updates cryptos list values per the Oct 4 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.646.61"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "1.644.50"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'213.06"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'0.531"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'22.96"
$ws.Range("E8").Value = "  -2.44%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "1.878.76"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "1.650.56"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "'4.04"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").Value = "'64.09"
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("D17").Value = "27.636.15"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").Value = "'228.99"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "'7.62"
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("D23").Value = "'10.04"
$ws.Range("E23").Value = "  +7.45%  "
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("D25").Value = "'148.91"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "'15.63"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("D30").Value = "'1.18"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").Value = "1.439.63"
$ws.Range("E34").Value = "  -0.55%  "
$ws.Range("E35").Value = "  +2.52%  "
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").Value = "'0.881"
$ws.Range("E38").Value = "  -2.72%  "
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").Value = "'0.893"
$ws.Range("E40").Value = "  +13.75%  "
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "'5.67"
$ws.Range("E43").Value = "  +3.66%  "
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("E45").Value = "  +1.99%  "
$ws.Range("D46").Value = "'65.26"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").Value = "1.787.65"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").Value = "'86.59"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0988"
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.72"
$ws.Range("E51").Value = "  +0.24%  "
